$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 10.46510533333333
$ws.Range("H2").Value = 31.395316
$ws.Range("I2").Value = 0.5554075997074865
$ws.Range("J2").Value = 0.5554075997074865
$ws.Range("M2").Value = 2.718682666666667
$ws.Range("N2").Value = 8.156048
$ws.Range("O2").Value = 0.5434637507613679
$ws.Range("P2").Value = 0.5434637507613679
$ws.Range("Q2").Value = 28.45130047457422
$ws.Range("R2").Value = 256.061704271168
$ws.Range("S2").Value = 0.3018438973383991
$ws.Range("T2").Value = 0.3018438973383991

# Row 3
$ws.Range("G3").Value = 10.46510533333333
$ws.Range("H3").Value = 31.395316
$ws.Range("I3").Value = 0.5554075997074865
$ws.Range("J3").Value = 0.5554075997074865
$ws.Range("M3").Value = 2.283827
$ws.Range("N3").Value = 6.851481
$ws.Range("O3").Value = 0.4565362492386322
$ws.Range("P3").Value = 0.4565362492386321
$ws.Range("Q3").Value = 23.90049011811067
$ws.Range("R3").Value = 215.104411062996
$ws.Range("S3").Value = 0.2535637023690875
$ws.Range("T3").Value = 0.2535637023690875

# Row 4
$ws.Range("I4").Value = 0.3053945925621632
$ws.Range("J4").Value = 0.3053945925621632
$ws.Range("M4").Value = 2.718682666666667
$ws.Range("N4").Value = 8.156048
$ws.Range("O4").Value = 0.5434637507613679
$ws.Range("P4").Value = 0.5434637507613679
$ws.Range("Q4").Value = 15.64413832448889
$ws.Range("R4").Value = 140.7972449204
$ws.Range("S4").Value = 0.1659708907360729
$ws.Range("T4").Value = 0.1659708907360729

# Row 5
$ws.Range("I5").Value = 0.3053945925621632
$ws.Range("J5").Value = 0.3053945925621632
$ws.Range("M5").Value = 2.283827
$ws.Range("N5").Value = 6.851481
$ws.Range("O5").Value = 0.4565362492386322
$ws.Range("P5").Value = 0.4565362492386321
$ws.Range("Q5").Value = 13.14184473799167
$ws.Range("R5").Value = 118.276602641925
$ws.Range("S5").Value = 0.1394237018260902
$ws.Range("T5").Value = 0.1394237018260902

# Row 6
$ws.Range("G6").Value = 2.146766
$ws.Range("H6").Value = 6.440298
$ws.Range("I6").Value = 0.1139338891693565
$ws.Range("J6").Value = 0.1139338891693565
$ws.Range("M6").Value = 2.718682666666667
$ws.Range("N6").Value = 8.156048
$ws.Range("O6").Value = 0.5434637507613679
$ws.Range("P6").Value = 0.5434637507613679
$ws.Range("Q6").Value = 5.836375513589333
$ws.Range("R6").Value = 52.527379622304
$ws.Range("S6").Value = 0.06191893874680848
$ws.Range("T6").Value = 0.0619189387468085

# Row 7
$ws.Range("G7").Value = 2.146766
$ws.Range("H7").Value = 6.440298
$ws.Range("I7").Value = 0.1139338891693565
$ws.Range("J7").Value = 0.1139338891693565
$ws.Range("M7").Value = 2.283827
$ws.Range("N7").Value = 6.851481
$ws.Range("O7").Value = 0.4565362492386322
$ws.Range("P7").Value = 0.4565362492386321
$ws.Range("Q7").Value = 4.902842153482
$ws.Range("R7").Value = 44.125579381338
$ws.Range("S7").Value = 0.05201495042254804
$ws.Range("T7").Value = 0.05201495042254804

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.4760280000000001
$ws.Range("H8").Value = 1.428084
$ws.Range("I8").Value = 0.02526391856099382
$ws.Range("J8").Value = 0.02526391856099382
$ws.Range("M8").Value = 2.718682666666667
$ws.Range("N8").Value = 8.156048
$ws.Range("O8").Value = 0.5434637507613679
$ws.Range("P8").Value = 0.5434637507613679
$ws.Range("Q8").Value = 1.294169072448
$ws.Range("R8").Value = 11.647521652032
$ws.Range("S8").Value = 0.01373002394008744
$ws.Range("T8").Value = 0.01373002394008744

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.4760280000000001
$ws.Range("H9").Value = 1.428084
$ws.Range("I9").Value = 0.02526391856099382
$ws.Range("J9").Value = 0.02526391856099382
$ws.Range("M9").Value = 2.283827
$ws.Range("N9").Value = 6.851481
$ws.Range("O9").Value = 0.4565362492386322
$ws.Range("P9").Value = 0.4565362492386321
$ws.Range("Q9").Value = 1.087165599156
$ws.Range("R9").Value = 9.784490392404001
$ws.Range("S9").Value = 0.01153389462090638
$ws.Range("T9").Value = 0.01153389462090638
